# fix(publipostage): Refactor synthetic array /3
#
# The "statut" color-code legend used throughout the sheet is renamed:
#   ⬛ (noir)  -> 📘 (bleu)
#   🟥         -> 📕
#   🟧         -> 📙
#   🟩         -> 📗
#
# These symbols/labels are repeated on every data row (columns A and B),
# so every matching cell across the sheet needs to be updated - not just
# a single occurrence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlWhole = 1

$all = $ws.Cells

$all.Replace("⬛", "📘", $xlWhole) | Out-Null
$all.Replace("🟥", "📕", $xlWhole) | Out-Null
$all.Replace("🟧", "📙", $xlWhole) | Out-Null
$all.Replace("🟩", "📗", $xlWhole) | Out-Null
$all.Replace("noir", "bleu", $xlWhole) | Out-Null
